# Daily attendance processing - 2025-10-23 13:34:33
# Update "Recorded By" (column G) values on the "Session Analysis Results" sheet
# so that the "System" entry is reordered within the comma-separated list of
# recorders for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$updates = @(
    @{ Row = 2; Value = 'system, System, backup@backdoor.com' },
    @{ Row = 3; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 4; Value = 'System, backup@backdoor.com' },
    @{ Row = 5; Value = 'System, backup@backdoor.com' },
    @{ Row = 6; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 8; Value = 'System, backup@backdoor.com' },
    @{ Row = 10; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 11; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 12; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 13; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 14; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 15; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 17; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 18; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 19; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 20; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 21; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 22; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 29; Value = 'system, System, backup@backdoor.com' },
    @{ Row = 30; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 31; Value = 'System, backup@backdoor.com' },
    @{ Row = 32; Value = 'System, backup@backdoor.com' },
    @{ Row = 33; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 35; Value = 'System, backup@backdoor.com' },
    @{ Row = 37; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 38; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 39; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 40; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 41; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 42; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 44; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 45; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 46; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 47; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 48; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 49; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 56; Value = 'system, System, backup@backdoor.com' },
    @{ Row = 57; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 58; Value = 'System, backup@backdoor.com' },
    @{ Row = 59; Value = 'System, backup@backdoor.com' },
    @{ Row = 60; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 62; Value = 'System, backup@backdoor.com' },
    @{ Row = 64; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 65; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 66; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 67; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 68; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 69; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 71; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 72; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 73; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 74; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 75; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 76; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 83; Value = 'System, backup@backdoor.com' },
    @{ Row = 84; Value = 'System, backup@backdoor.com' },
    @{ Row = 85; Value = 'System, backup@backdoor.com' },
    @{ Row = 86; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 87; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 88; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 89; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 93; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 95; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 96; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 97; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 99; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 102; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 109; Value = 'System, backup@backdoor.com' },
    @{ Row = 110; Value = 'System, backup@backdoor.com' },
    @{ Row = 111; Value = 'System, backup@backdoor.com' },
    @{ Row = 112; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 113; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 114; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 115; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 119; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 121; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 122; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 123; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 125; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 128; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 135; Value = 'System, backup@backdoor.com' },
    @{ Row = 136; Value = 'System, backup@backdoor.com' },
    @{ Row = 137; Value = 'System, backup@backdoor.com' },
    @{ Row = 138; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 139; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 140; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 141; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 145; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 147; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 148; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 149; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 151; Value = 'System, dnasr281@gmail.com' },
    @{ Row = 154; Value = 'System, dnasr281@gmail.com' }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Value
}
